$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row 11 (shifts old rows 11-17 down to 12-18)
$ws.Rows.Item(11).Insert()

# Populate the new row with the waittopageload1 step + its amount
$ws.Range("A11").Value = "waittopageload1"
$ws.Range("B11").Value = 2000

# Match the number formatting used by the analogous "amount" row (B3)
$ws.Range("B3").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select the newly inserted row's range, then activate this sheet so it
# becomes the active tab (tabSelected / workbookView activeTab follow this)
[void]$ws.Range("A11:B11").Select()
[void]$ws.Activate()
